$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D stores prices as plain text that can look numeric (e.g. "1.002",
# "29.925.49"). Prefix with a leading apostrophe (quote-prefix) so Excel
# keeps these as literal text instead of coercing them into numbers.

$ws.Range("D2").Value = "'29.925.49"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "'1.903.88"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'318.98"
$ws.Range("E5").Value = "  -2.35%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D7").Value = "'0.5021"
$ws.Range("E7").Value = "  -3.04%  "
$ws.Range("D8").Value = "'0.4035"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.08240"
$ws.Range("E9").Value = "  -3.05%  "
$ws.Range("D10").Value = "'41.96"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").Value = "'1.095"
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("D12").Value = "'23.91"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").Value = "'1.904.78"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "'6.348"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "'7.179"
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "'91.74"
$ws.Range("E17").Value = "  -3.68%  "
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").Value = "'0.06490"
$ws.Range("E19").Value = "  -2.80%  "
$ws.Range("D20").Value = "'17.97"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'5.922"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").Value = "'29.954.73"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "'2.195"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").Value = "'22.04"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("D27").Value = "'2.127.69"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "'161.39"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "'2.254"
$ws.Range("E29").Value = "  -6.60%  "
$ws.Range("D30").Value = "'128.47"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").Value = "'1.119"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").Value = "'5.897"
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("D34").Value = "'3.792"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").Value = "'5.364"
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("D36").Value = "'0.02429"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D37").Value = "'0.06327"
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("D38").Value = "'0.2134"
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.188"
$ws.Range("E39").Value = "  -3.87%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6436"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").Value = "'8.614"
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("D42").Value = "'11.26"
$ws.Range("E42").Value = "  -5.65%  "
$ws.Range("D43").Value = "'1.200"
$ws.Range("E43").Value = "  -3.34%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.28"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'2.179"
$ws.Range("E45").Value = "  +5.25%  "
$ws.Range("D46").Value = "'0.5993"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").Value = "'122.15"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("D49").Value = "'1.201"
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("D50").Value = "'78.20"
$ws.Range("E50").Value = "  -1.82%  "
$ws.Range("D51").Value = "'1.129"
